$wb = $excel.ActiveWorkbook

# --- Add the new "Empty" worksheet, placed right after "trim-range" ---
$trimRange = $wb.Worksheets.Item("trim-range")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $trimRange)
$newSheet.Name = "Empty"

# Renaming a sheet causes the engine to rewrite every defined name and
# drops the sheet-qualifier on the "Missing" name's #REF! error. Restore
# its original RefersTo text so definedNames stays byte-for-byte the same.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Missing") {
        $n.RefersTo = "='trim-range'!#REF!"
    }
}

# --- Populate the new sheet ---
# A1: literal text "=C1" (quote-prefixed, not a live formula)
$newSheet.Range("A1").Value = "'=C1"
# B1: real formula referencing the (empty) C1 -> evaluates to 0
$newSheet.Range("B1").Formula = "=C1"

# --- Selections / active-tab bookkeeping ---
# trim-range keeps cell G1 selected but is no longer the active tab.
$trimRange.Range("G1").Select() | Out-Null
# Empty becomes the active sheet/tab with C1 selected.
$newSheet.Activate() | Out-Null
$newSheet.Range("C1").Select() | Out-Null
